$d = $word.ActiveDocument
$d.Content.Find.Execute("Le devoir s’oppose-t-il à la liberté ?", $true, $false, $false, $false, $false, $true, 1, $false, "Le devoir s’oppose-t-il à la liberté ? πα (pentaraïe)", 2)
